# Weekly update: a new price-report row is inserted above the former row 48
# (everything from the old row 48 down to row 137 shifts down by one row,
# ending at row 138), and the newly opened row 48 is populated with the
# latest "Crespo record" / "Primera" report for Terminal Hortofrutícola
# Agro Chillán.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 48:137 down to 49:138, leaving row 48 free for the new record.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new weekly record.
$ws.Cells.Item(48, 1).Value = 7
$ws.Cells.Item(48, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(48, 3).Value = "Ñuble"
$ws.Cells.Item(48, 4).Value = 44495
$ws.Cells.Item(48, 5).Value = 16
$ws.Cells.Item(48, 6).Value = 100112006
$ws.Cells.Item(48, 7).Value = "Repollo"
$ws.Cells.Item(48, 8).Value = "Crespo record"
$ws.Cells.Item(48, 9).Value = "Primera"
$ws.Cells.Item(48, 10).Value = 180
$ws.Cells.Item(48, 11).Value = 600
$ws.Cells.Item(48, 12).Value = 700
$ws.Cells.Item(48, 13).Value = 650
$ws.Cells.Item(48, 14).Value = "$/unidad"
$ws.Cells.Item(48, 15).Value = "Región del Maule"
$ws.Cells.Item(48, 16).Value = 650
$ws.Cells.Item(48, 17).Value = 1
$ws.Cells.Item(48, 18).Value = "Hortaliza"

# Keep the same date number format the rest of column D uses.
$ws.Cells.Item(48, 4).NumberFormat = $ws.Cells.Item(49, 4).NumberFormat
